# Adds the "ODI Bowling Extra" sheet (MATCH_CODE / MAIDEN_OVERS /
# PERCENT_WICKETS_OF_ALL) as the 5th, trailing sheet of the workbook, mirroring
# the layout/format already used by the other "*_Extra" sheets (e.g.
# "ODI Batting Extra").

$wb = $excel.ActiveWorkbook

# --- create the new sheet after the last existing sheet -------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([Type]::Missing, $lastSheet)
$ws.Name = "ODI Bowling Extra"

# --- header row -------------------------------------------------------------
$headers = @("MATCH_CODE", "MAIDEN_OVERS", "PERCENT_WICKETS_OF_ALL")
for ($col = 1; $col -le $headers.Length; $col++) {
    $cell = $ws.Cells.Item(1, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $headers[$col - 1]
}

$headerRange = $ws.Range("A1:C1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108   # xlCenter
$headerRange.VerticalAlignment = -4160     # xlTop
$headerRange.Borders.LineStyle = 1         # xlContinuous

# --- data rows ---------------------------------------------------------------
# MATCH_CODE, MAIDEN_OVERS, PERCENT_WICKETS_OF_ALL
$data = @(
    @('3834', '0', '10.00%'),
    @('4069', '0', '10.00%'),
    @('4071', '0', '10.00%'),
    @('4166', '0', ''),
    @('4167', '', ''),
    @('4168', '0', '10.00%'),
    @('4169', '', ''),
    @('4170', '2', '10.00%'),
    @('4398', '1', '10.00%'),
    @('4399', '0', ''),
    @('4400', '0', '10.00%'),
    @('4419', '', ''),
    @('4437', '0', '20.00%'),
    @('4486', '', ''),
    @('4594', '', ''),
    @('4645', '', ''),
    @('4646', '0', '10.00%'),
    @('4660', '', ''),
    @('4663', '', ''),
    @('4732', '0', '20.00%')
)

$rowIndex = 2
foreach ($record in $data) {
    $ws.Cells.Item($rowIndex, 1).NumberFormat = "@"
    $ws.Cells.Item($rowIndex, 1).Value = $record[0]

    $ws.Cells.Item($rowIndex, 2).NumberFormat = "@"
    if ($record[1] -ne '') {
        $ws.Cells.Item($rowIndex, 2).Value = $record[1]
    }

    $ws.Cells.Item($rowIndex, 3).NumberFormat = "@"
    if ($record[2] -ne '') {
        $ws.Cells.Item($rowIndex, 3).Value = $record[2]
    }

    $rowIndex++
}
